$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) for another shipment entry, delivering the cod
# count directly via the existing path-formula columns instead of the
# removed cod_count_per_user helper.
$ws.Range("A3").Value = "جيانا"
$ws.Range("B3").Formula = "=""samples\""&A3&"" ""&C3"
$ws.Range("C3").Value = "02-02"
$ws.Range("D3").Value = "1"

$ws.Range("B6").Select() | Out-Null
